$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Trace3")
$ws.Range("T2").Value = 8504660.8790714014

$ws = $wb.Worksheets.Item("Trace4")
$ws.Range("T2").Value = 8482088.6430714

$ws = $wb.Worksheets.Item("Trace5")
$ws.Range("T2").Value = 8496735.4630714022

$ws = $wb.Worksheets.Item("Trace6")
$ws.Range("T2").Value = 8511236.7730714008

$ws = $wb.Worksheets.Item("Trace7")
$ws.Range("T2").Value = 8551216.5370714013

$ws = $wb.Worksheets.Item("Trace8")
$ws.Range("T2").Value = 8477104.2890713997

$ws = $wb.Worksheets.Item("Trace9")
$ws.Range("T2").Value = 8461287.2070714012

$ws = $wb.Worksheets.Item("Trace10")
$ws.Range("T2").Value = 8505758.7370714005

$ws = $wb.Worksheets.Item("Trace11")
$ws.Range("T2").Value = 8472915.5890714005
$ws.Range("R3").Value = 8622692.9817537088
$ws.Range("T3").Value = 8773290.6517537106

$ws = $wb.Worksheets.Item("Trace12")
$ws.Range("T2").Value = 8490785.967071401

$ws = $wb.Worksheets.Item("Trace13")
$ws.Range("T2").Value = 8500876.0090714004

$ws = $wb.Worksheets.Item("Trace14")
$ws.Range("T2").Value = 8495807.2550713997

$ws = $wb.Worksheets.Item("Trace15")
$ws.Range("T2").Value = 8501376.8970714007

$ws = $wb.Worksheets.Item("Trace16")
$ws.Range("T2").Value = 8471722.4650714006
$ws.Range("R3").Value = 8999999.9999775533
$ws.Range("T3").Value = 9132878.8059775531

$ws = $wb.Worksheets.Item("Trace17")
$ws.Range("T2").Value = 8462468.2210714016

$ws = $wb.Worksheets.Item("Trace18")
$ws.Range("T2").Value = 8502177.047071401

$ws = $wb.Worksheets.Item("Trace19")
$ws.Range("T2").Value = 8472836.3930714

$ws = $wb.Worksheets.Item("Trace20")
$ws.Range("T2").Value = 8443871.3290714025

$ws = $wb.Worksheets.Item("Trace21")
$ws.Range("T2").Value = 8420901.9490714017

$ws = $wb.Worksheets.Item("Trace22")
$ws.Range("T2").Value = 8499657.1190714017

$ws = $wb.Worksheets.Item("Trace23")
$ws.Range("T2").Value = 8506393.9590714015

$ws = $wb.Worksheets.Item("Trace24")
$ws.Range("T2").Value = 8456777.5090714004

$ws = $wb.Worksheets.Item("Trace25")
$ws.Range("T2").Value = 8487124.3090713993

$ws = $wb.Worksheets.Item("Trace26")
$ws.Range("T2").Value = 8501401.8390714023

$ws = $wb.Worksheets.Item("Trace27")
$ws.Range("T2").Value = 8457165.4790714029

$ws = $wb.Worksheets.Item("Trace28")
$ws.Range("T2").Value = 8501285.0290714018

$ws = $wb.Worksheets.Item("Trace29")
$ws.Range("T2").Value = 8546595.8790713996

$ws = $wb.Worksheets.Item("Trace30")
$ws.Range("T2").Value = 8480828.2490714006

$ws = $wb.Worksheets.Item("Trace31")
$ws.Range("T2").Value = 8508802.8090713993

$ws = $wb.Worksheets.Item("Trace32")
$ws.Range("T2").Value = 8502448.399071401

$ws = $wb.Worksheets.Item("Trace33")
$ws.Range("T2").Value = 8494660.4290714003

$ws = $wb.Worksheets.Item("Trace34")
$ws.Range("T2").Value = 8519276.1590714008

$ws = $wb.Worksheets.Item("Trace35")
$ws.Range("T2").Value = 8446422.5990714021

$ws = $wb.Worksheets.Item("Trace36")
$ws.Range("T2").Value = 8471200.6890714001

$ws = $wb.Worksheets.Item("Trace37")
$ws.Range("T2").Value = 8486210.149071401

$ws = $wb.Worksheets.Item("Trace38")
$ws.Range("T2").Value = 8484151.0690714009
